# Apply the StructureDefinition-employee-count.xlsx update:
#  - bump URL/Version/Date/Publisher on the "Metadata" sheet
#  - clear the stray Constraint(s) text on the base "Extension" row of
#    the "Elements" sheet (row 2, column AI) - it only belongs on the
#    Extension.extension row (row 4), not the top-level Extension row.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-count"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
